$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Range("B18").Value = "[Leandro-M. Maq. E. I., Leandro-M.S.R.A.C., Guilherme-Coman. Hidraulicos, Victor-Usin. CNC]"
$ws.Range("C18").Value = "[Cláudio-Soldagem, Ismail-Metrologia 2, Leandro-M.S.R.A.C., Guilherme-C. L. P.]"
$ws.Range("D18").Value = "Suzanny-Des. Maq. "
$ws.Range("E18").Value = "[Paulo Rob.-CAD/CAM, Guilherme-C. L. P., Guilherme-Coman. Hidraulicos, Aderci-Fresagem]"
$ws.Range("F18").Value = "Claudinei-Elem"

# Row 19
$ws.Range("B19").Value = "[Leandro-M. Maq. E. I., Leandro-M.S.R.A.C., Joel L.-Fundição, Victor-Usin. CNC]"
$ws.Range("C19").Value = "[Cláudio-Soldagem, Ismail-Metrologia 2, Leandro-M. Maq. E. I., Rogério-Retífica]"
$ws.Range("E19").Value = "[Paulo Rob.-CAD/CAM, Guilherme-C. L. P., Guilherme-Coman. Hidraulicos, Aderci-Fresagem]"
$ws.Range("F19").Value = "Claudinei-Elem"

# Row 20
$ws.Range("B20").Value = "[Leandro-M. Maq. E. I., Leandro-M.S.R.A.C., Joel L.-Fundição, Elcio Dec.-C.pneumática]"
$ws.Range("C20").Value = "[Cláudio-Soldagem, Ismail-Metrologia 2, Elcio Dec.-C.pneumática, Rogério-Retífica]"
$ws.Range("D20").Value = "Suzanny-Des. Maq. "
$ws.Range("E20").Value = "[Paulo Rob.-CAD/CAM, Ismail-Metrologia 2, Cláudio-Soldagem, Aderci-Fresagem]"
$ws.Range("F20").Value = "Euclides-Gestão integr"

# Row 21
$ws.Range("B21").Value = "[Elcio Dec.-C.pneumática, Rogério-Retífica, Joel L.-Fundição, Victor-Usin. CNC]"
$ws.Range("C21").Value = "[-, Joel L.-Fundição, Elcio Dec.-C.pneumática, Rogério-Retífica]"
$ws.Range("D21").Value = "[-, -, Victor-Usin. CNC, -]"
$ws.Range("E21").Value = "[Guilherme-C. L. P., Paulo Rob.-CAD/CAM, Guilherme-Coman. Hidraulicos, Aderci-Fresagem]"
$ws.Range("F21").Value = "Euclides-Gestão integr"
